# Re-generate the linear/quadratic problem data (matching commit:
# "volver a generar problemas cuadraticos y lineales").
#
# This rewrites the numeric/expression cells on several sheets with a new
# set of randomly-generated coefficients. Cells that hold numbers stored
# as text (e.g. "7.1") must stay text-typed, so we force NumberFormat="@"
# before writing them (this mirrors what happens in Excel's UI when typing
# into a cell that was already formatted as Text).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Restricciones_del_follower : Expression / Function_Evaluation /
# Restriction_Set_Type / Lambda_value / Beta_value / Gamma_value
# ---------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

$followerRange = $wsFollower.Range("A2:F6")
$followerRange.NumberFormat = "@"

$wsFollower.Range("A2").Value = "3.728814523931479 - 0.8427522622502988y_1 + 1.3278697854419215y_2"
$wsFollower.Range("B2").Value = "-3.728814523931479"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.75"
$wsFollower.Range("E2").Value = "9.5"
$wsFollower.Range("F2").Value = "2.8000000000000003"

$wsFollower.Range("A3").Value = "-8.17158955096466 + 0.9861703943998634y_1 - 0.11678333617893122y_2"
$wsFollower.Range("B3").Value = "4.17158955096466"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.19"
$wsFollower.Range("E3").Value = "5.5"
$wsFollower.Range("F3").Value = "0"

$wsFollower.Range("A4").Value = "29.207453846923947 - 2x - 0.10425783580457493y_1 - 5.324843946794189y_2"
$wsFollower.Range("B4").Value = "-45.20745384692395"
$wsFollower.Range("C4").Value = "J_0_LP_v"
$wsFollower.Range("D4").Value = "0.24"
$wsFollower.Range("E4").Value = "0"
$wsFollower.Range("F4").Value = "1.2"

$wsFollower.Range("A5").Value = "-72.83043423823345 + 8x + 1.2217858972170053y_1 + 1.872858687610267y_2"
$wsFollower.Range("B5").Value = "24.270434238233452"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.92"
$wsFollower.Range("E5").Value = "0"
$wsFollower.Range("F5").Value = "4.6000000000000005"

$wsFollower.Range("A6").Value = "-6.543179101929319 - 2x - 1.9723407887997269y_1 + 0.23356667235786244y_2"
$wsFollower.Range("B6").Value = "-18.54317910192932"
$wsFollower.Range("C6").Value = "J_Ne_L0_v"
$wsFollower.Range("D6").Value = "0.82"
$wsFollower.Range("E6").Value = "3.4000000000000004"
$wsFollower.Range("F6").Value = "0"

# ---------------------------------------------------------------------
# Punto_modificado : x / y_1 / y_2
# ---------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
$wsPunto.Range("A2:C2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "7.1"
$wsPunto.Range("B2").Value = "8.600000000000001"
$wsPunto.Range("C2").Value = "2.65"

# ---------------------------------------------------------------------
# Vector_bf
# ---------------------------------------------------------------------
$wsBf = $wb.Worksheets.Item("Vector_bf")
$wsBf.Range("A2:A3").NumberFormat = "@"
$wsBf.Range("A2").Value = "1.9629901237209788"
$wsBf.Range("A3").Value = "-1.6103056219117318"

# ---------------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------------
$wsBF = $wb.Worksheets.Item("Vector_BF")
$wsBF.Range("A2:A4").NumberFormat = "@"
$wsBF.Range("A2").Value = "7.800000000000001"
$wsBF.Range("A3").Value = "12.28816800409766"
$wsBF.Range("A4").Value = "-14.766581298730864"

# ---------------------------------------------------------------------
# Vector_Alpha (stored as real numbers, not text)
# ---------------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 0.27
$wsAlpha.Range("A3").Value = 2.2800000000000002

Write-Host "Regenerated linear/quadratic problem data."
